$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARN: could not find [$old]"
    }
    return $ok
}

# --- Title ---
Replace-Text "Unveiling the Enigma of Consciousness" "Delving into the World of Biology: From Atoms to Ecosystems" | Out-Null

# --- Author byline ---
Replace-Text " Kristina Petrova" " Sarah LeBlanc" | Out-Null

# --- Email paragraph (run-level edits, then append new runs) ---
Replace-Text "kpetrova@sciencereview" "dr" | Out-Null
Replace-Text "org" "sarahleblanc@growthandlife" | Out-Null

$emailPara = $d.Paragraphs.Item(3)
$endOfEmail = $d.Range($emailPara.Range.End - 1, $emailPara.Range.End - 1)
$endOfEmail.InsertAfter(".")
$emailPara2 = $d.Paragraphs.Item(3)
$endOfEmail2 = $d.Range($emailPara2.Range.End - 1, $emailPara2.Range.End - 1)
$endOfEmail2.InsertAfter("edu")

Write-Output "Paragraph 3 now:"
Write-Output $d.Paragraphs.Item(3).Range.Text

# --- Body paragraph (paragraph 5) ---
Replace-Text "From ancient mystics to modern neuroscientists, the nature of consciousness has intrigued thinkers across disciplines" "Biology, an intricate symphony of life, unveils the captivating saga of the cosmos, from the tiniest particles to the boundless expanse of biosystems" | Out-Null

Replace-Text " Consciousness, the subjective experience of our internal and external worlds, remains one of the greatest unsolved mysteries in science today" " In this exploration, we delve into the realm of living things, unraveling the secrets of cells, metabolism, and evolution" | Out-Null

Replace-Text " How does the physical substrate of our brain produce the rich tapestry of subjective experience? Why do qualities like awareness, intentionality, and introspection arise from neural processes? Delving into the enigma of consciousness, we explore the boundary between physical reality and subjective experience" " Intriguing questions arise: How do organisms function? How do they adapt to their surroundings? How are life's diverse patterns interconnected? From the smallest interactions within cells to the grand tapestry of ecosystems, biology weaves a compelling narrative of unity and diversity in the living world" | Out-Null

Replace-Text "Our investigation begins by examining the complex interplay between brain structure and function" "Advancing further into this realm, we encounter the fascinating realm of genetics, where the thread of life is encoded" | Out-Null

Replace-Text " Neuroscientists have identified specific brain regions involved in conscious processing, including the cerebral cortex, thalamus, and brainstem" " Decoding the language of DNA reveals the blueprint for heredity, variation, and the evolution of species" | Out-Null

Replace-Text " However, the exact mechanisms by which these neural networks generate consciousness remain elusive" " Witnessing the dynamism of evolution, we unravel intricate threads of change, shaping organisms and ecosystems over eons" | Out-Null

Replace-Text " We are confronted with the hard problem of consciousness - the challenge of explaining how physical processes manifest subjective experiences" " Embracing a dynamic view of life, we navigate the intricate interplay of systems within organisms, exploring how their coordinated interactions maintain a delicate balance" | Out-Null

Replace-Text "Further complicating our understanding, the concept of consciousness varies across cultures and philosophical traditions" "Highlighting the pervasive presence of biology in our lives, medical discoveries and biotechnological wonders unfold before our eyes" | Out-Null

Replace-Text " Some assert that consciousness is an emergent property of complex systems, arising from the interactions of billions of neurons" " From the development of vaccines and treatments to exploring gene editing techniques, biology empowers us to safeguard health and unravel the secrets of longevity" | Out-Null

Replace-Text " Others posit that consciousness is non-physical, existing outside the realm of space and time" " Agriculture and conservation benefit from our understanding of biological principles, ensuring food security and safeguarding biodiversity" | Out-Null

Replace-Text " The mystery deepens as we encounter altered states of consciousness, such as dreams, meditation, and psychedelic experiences, that challenge our conventional notions of conscious awareness" " The journey of biology is a transformative exploration, unlocking the secrets of life and inspiring us to cherish and sustain the wonders of our living world" | Out-Null

Write-Output "Paragraph 5 now:"
Write-Output $d.Paragraphs.Item(5).Range.Text

# --- Summary paragraph (paragraph 7) ---
Replace-Text "Our exploration of consciousness reveals a complex phenomenon that defies easy explanation" "Biology captivates us with its exploration of living systems, from the fundamental principles of cells to the interconnectedness of ecosystems" | Out-Null

# " The relationship ... disciplines" -> " Genetics unveils the secrets of heredity " + new run (lastRenderedPageBreak) "and evolution, ... well-being"
Replace-Text " The relationship between brain activity and subjective experience remains an enigma, inviting further inquiry and innovation across scientific disciplines" " Genetics unveils the secrets of heredity " | Out-Null

$summaryPara = $d.Paragraphs.Item(7)
$fullText = $summaryPara.Range.Text
$marker = " Genetics unveils the secrets of heredity "
$insertPos = $summaryPara.Range.Start + $fullText.IndexOf($marker) + $marker.Length
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter("and evolution, while medical and technological advancements highlight biology's immense impact on human well-being")

Write-Output "Paragraph 7 now (after genetics split):"
Write-Output $d.Paragraphs.Item(7).Range.Text

# " As we continue ... universe" -> " Agriculture and conservation leverage ... future" + new run "." + new run " Immersed in the wonders ... ecosystems"
Replace-Text " As we continue to unravel the mysteries of consciousness, we may gain insights into the profound nature of our existence, shedding light on the essence of being human and our place in the universe" " Agriculture and conservation leverage biological understanding to ensure a sustainable and flourishing future" | Out-Null

$summaryPara2 = $d.Paragraphs.Item(7)
$fullText2 = $summaryPara2.Range.Text
$marker2 = " Agriculture and conservation leverage biological understanding to ensure a sustainable and flourishing future"
$insertPos2 = $summaryPara2.Range.Start + $fullText2.IndexOf($marker2) + $marker2.Length
$insertRange2 = $d.Range($insertPos2, $insertPos2)
$insertRange2.InsertAfter(".")
$insertPos3 = $insertPos2 + 1
$insertRange3 = $d.Range($insertPos3, $insertPos3)
$insertRange3.InsertAfter(" Immersed in the wonders of biology, we gain a deeper appreciation for the intricate dance of life, inspiring us to protect and nurture the delicate balance of our planet's ecosystems")

Write-Output "Paragraph 7 final:"
Write-Output $d.Paragraphs.Item(7).Range.Text
